# Update "想去人数" (want-to-go count) figures for several events across
# sheets "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value2 = 3070
$wsExpo.Range("F4").Value2 = 110
$wsExpo.Range("F5").Value2 = 6804
$wsExpo.Range("F6").Value2 = 1832
$wsExpo.Range("F7").Value2 = 51
$wsExpo.Range("F9").Value2 = 33
$wsExpo.Range("F11").Value2 = 136
$wsExpo.Range("F12").Value2 = 154

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value2 = 3

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value2 = 3070
$wsAll.Range("F3").Value2 = 3
$wsAll.Range("F5").Value2 = 110
$wsAll.Range("F6").Value2 = 6804
$wsAll.Range("F7").Value2 = 1832
$wsAll.Range("F8").Value2 = 51
$wsAll.Range("F10").Value2 = 33
$wsAll.Range("F12").Value2 = 136
$wsAll.Range("F13").Value2 = 154
